# Refresh the workbook with the latest published NHS/ONS COVID-19 deaths-by-ethnicity data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COVID_deaths_by_ethnicity")
$ws.Range("C3").Value = "Data for England up until 19th May 2020 (published 21st May 2020)"

$raw = $wb.Worksheets.Item("raw_data")

# Row 2: Any other Asian background
$raw.Range("D2").Value = 0.01545678
$raw.Range("E2").Value = 374
$raw.Range("F2").Value = 0.01643089
$raw.Range("G2").Value = 352
$raw.Range("H2").Value = 110
$raw.Range("I2").Value = 0.00482198
$raw.Range("J2").Value = 264
$raw.Range("K2").Value = 3.4

# Row 3: Bangladeshi
$raw.Range("D3").Value = 0.00823418
$raw.Range("E3").Value = 152
$raw.Range("F3").Value = 0.0066778
$raw.Range("G3").Value = 187
$raw.Range("H3").Value = 39
$raw.Range("I3").Value = 0.00169406
$raw.Range("J3").Value = 113
$raw.Range("K3").Value = 3.8974359

# Row 4: Chinese
$raw.Range("D4").Value = 0.00715875
$raw.Range("E4").Value = 78
$raw.Range("F4").Value = 0.00342676
$raw.Range("G4").Value = 163
$raw.Range("H4").Value = 56
$raw.Range("I4").Value = 0.00245445
$raw.Range("J4").Value = 22
$raw.Range("K4").Value = 1.39285714

# Row 5: Indian
$raw.Range("D5").Value = 0.02632781
$raw.Range("E5").Value = 695
$raw.Range("F5").Value = 0.03053335
$raw.Range("G5").Value = 599
$raw.Range("H5").Value = 290
$raw.Range("I5").Value = 0.0127392
$raw.Range("J5").Value = 405
$raw.Range("K5").Value = 2.39655172

# Row 6: Pakistani
$raw.Range("D6").Value = 0.02098152
$raw.Range("E6").Value = 462
$raw.Range("F6").Value = 0.02029699
$raw.Range("G6").Value = 478
$raw.Range("H6").Value = 127
$raw.Range("I6").Value = 0.00556962
$raw.Range("J6").Value = 335
$raw.Range("K6").Value = 3.63779528

# Row 7: African
$raw.Range("D7").Value = 0.01844361
$raw.Range("E7").Value = 400
$raw.Range("F7").Value = 0.01757315
$raw.Range("G7").Value = 420
$raw.Range("H7").Value = 78
$raw.Range("I7").Value = 0.00343592
$raw.Range("J7").Value = 322
$raw.Range("K7").Value = 5.12820513

# Row 8: Any other Black background
$raw.Range("D8").Value = 0.00524135
$raw.Range("E8").Value = 206
$raw.Range("F8").Value = 0.00905017
$raw.Range("G8").Value = 119
$raw.Range("H8").Value = 26
$raw.Range("I8").Value = 0.00113204
$raw.Range("J8").Value = 180
$raw.Range("K8").Value = 7.92307692

# Row 9: Caribbean
$raw.Range("D9").Value = 0.01114863
$raw.Range("E9").Value = 609
$raw.Range("F9").Value = 0.02675512
$raw.Range("G9").Value = 254
$raw.Range("H9").Value = 191
$raw.Range("I9").Value = 0.0084027
$raw.Range("J9").Value = 418
$raw.Range("K9").Value = 3.18848168

# Row 10: Any other Mixed background
$raw.Range("D10").Value = 0.00533846
$raw.Range("E10").Value = 68
$raw.Range("F10").Value = 0.00298744
$raw.Range("G10").Value = 122
$raw.Range("H10").Value = 27
$raw.Range("I10").Value = 0.00117646
$raw.Range("J10").Value = 41
$raw.Range("K10").Value = 2.51851852

# Row 11: White and Asian
$raw.Range("D11").Value = 0.00627603
$raw.Range("E11").Value = 29
$raw.Range("F11").Value = 0.00127405
$raw.Range("G11").Value = 143
$raw.Range("H11").Value = 28
$raw.Range("I11").Value = 0.00121093
$raw.Range("J11").Value = 1
$raw.Range("K11").Value = 1.03571429

# Row 12: White and Black African
$raw.Range("D12").Value = 0.0030474
$raw.Range("E12").Value = 14
$raw.Range("F12").Value = 0.000615
$raw.Range("G12").Value = 69
$raw.Range("H12").Value = 10
$raw.Range("I12").Value = 0.000428
$raw.Range("J12").Value = 4
$raw.Range("K12").Value = 1.4

# Row 13: White and Black Caribbean
$raw.Range("D13").Value = 0.00783997
$raw.Range("E13").Value = 45
$raw.Range("F13").Value = 0.00197698
$raw.Range("G13").Value = 178
$raw.Range("H13").Value = 37
$raw.Range("I13").Value = 0.00162175
$raw.Range("J13").Value = 8
$raw.Range("K13").Value = 1.21621622

# Row 14: Any other ethnic group
$raw.Range("D14").Value = 0.01034508
$raw.Range("E14").Value = 608
$raw.Range("F14").Value = 0.02671119
$raw.Range("G14").Value = 235
$raw.Range("H14").Value = 75
$raw.Range("I14").Value = 0.00327653
$raw.Range("J14").Value = 533
$raw.Range("K14").Value = 8.10666667

# Row 15: Any other White background
$raw.Range("D15").Value = 0.04687398
$raw.Range("E15").Value = 780
$raw.Range("F15").Value = 0.03426764
$raw.Range("G15").Value = 1067
$raw.Range("H15").Value = 466
$raw.Range("I15").Value = 0.02046486
$raw.Range("J15").Value = 314
$raw.Range("K15").Value = 1.67381974

# Row 16: British
$raw.Range("D16").Value = 0.797534
$raw.Range("E16").Value = 18019
$raw.Range("F16").Value = 0.79162639
$raw.Range("G16").Value = 18153
$raw.Range("H16").Value = 20837
$raw.Range("I16").Value = 0.91542614
$raw.Range("J16").Value = -2818
$raw.Range("K16").Value = 0.8647598

# Row 17: Irish
$raw.Range("D17").Value = 0.00975244
$raw.Range("E17").Value = 223
$raw.Range("F17").Value = 0.00979703
$raw.Range("G17").Value = 222
$raw.Range("H17").Value = 368
$raw.Range("I17").Value = 0.01614573
$raw.Range("J17").Value = -145
$raw.Range("K17").Value = 0.60597826

# Restore the active selection on the summary sheet to match the saved view.
$ws.Activate()
$ws.Range("P22").Select()

